$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header column AH1 -> "convivenza" (inherits bold/centered header style
# from the header row's row-level formatting, matching AG1..A1)
$ws.Range("AH1").Value = "convivenza"

# New data row 17 (questionario docenti + studenti upload)

# Plain text cells (letters or mixed punctuation -> Excel keeps these as text
# automatically).
$ws.Range("A17").Value = "2024-11-27 10:27:13"
$ws.Range("B17").Value = "2024-11-27 10:53:38"
$ws.Range("D17").Value = "109.112.95.6"
$ws.Range("I17").Value = "R_2eWC9Qx590lQ5ae"
$ws.Range("P17").Value = "anonymous"
$ws.Range("Q17").Value = "IT"
$ws.Range("Z17").Value = "Classe abbastanza integrata e rispettosa delle regole"
$ws.Range("AA17").Value = "Buone relazioni nel gruppo classe"
$ws.Range("AB17").Value = "Buona partecipazione"
$ws.Range("AC17").Value = "Buono il rispetto degli impegni scolastici"
$ws.Range("AE17").Value = "buona collaborazione"
$ws.Range("AF17").Value = "E' una classe molto curiosa"
$ws.Range("AG17").Value = "BR03"
$ws.Range("AH17").Value = "Buono"

# Cells whose literal text looks like a pure number ("0", "100", "1585", a
# coordinate, etc). The source file stores these as text (inlineStr), same
# as the matching columns in every other data row, so force a Text number
# format while typing them in, then clear the format again so the cell
# doesn't end up carrying a stray style index.
$textCells  = @("C17","E17","F17","G17","H17","N17","O17","R17","S17")
$textValues = @("0","100","1585","1","1732704819.245","41.8904","12.5126","1","1")
for ($i = 0; $i -lt $textCells.Length; $i++) {
    $ws.Range($textCells[$i]).NumberFormat = "@"
    $ws.Range($textCells[$i]).Value = $textValues[$i]
    $ws.Range($textCells[$i]).ClearFormats()
}

# True numeric cells (stored as <v> without a text type in the target).
$ws.Range("T17").Value = 27
$ws.Range("U17").Value = 14
$ws.Range("V17").Value = 13
$ws.Range("W17").Value = 11
$ws.Range("X17").Value = 7
$ws.Range("Y17").Value = 4
$ws.Range("AD17").Value = 4
